$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1615418222.678129
$ws.Range("C3").Value = 1615418233.622086
$ws.Range("C4").Value = 1615418236.867699
$ws.Range("C5").Value = 1615418486.301492
$ws.Range("C6").Value = 1615418486.43469
$ws.Range("C7").Value = 1615418493.962464
$ws.Range("C8").Value = 1615418493.962464
